# "Trabajando con cambios y errores en el sistema"
# Fill in missing "Primer Apellido" / "Nombre(s)" data that had been left
# blank for a handful of beneficiarios, and fix a bad "Integrantes familia"
# count on row 6. Also update the saved view (active cell / scroll
# position) to match where the author left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (AAAA730807MZSRRG00): Primer Apellido was blank -> same apellido as
# "Segundo Apellido" (ARAUJO).
$ws.Range("B3").Value = "ARAUJO"

# Row 4 (AAAA830602MZSRVL02): Primer Apellido was blank -> ARANDA.
$ws.Range("B4").Value = "ARANDA"

# Row 5 (AAAB940222MZSLLR02): Nombre(s) was blank -> BRENDA ALEJANDRA.
$ws.Range("D5").Value = "BRENDA ALEJANDRA"

# Row 6 (AAAB960415MZSLNR09): Nombre(s) was blank -> BRENDA PAOLA, and the
# "Integrantes familia" value was wrong (8 -> 1).
$ws.Range("D6").Value = "BRENDA PAOLA"
$ws.Range("E6").Value = 1

# Row 8 (AAAD710729HZSNGV09): Nombre(s) was blank -> DAVID.
$ws.Range("D8").Value = "DAVID"

# Update the view state left in the file: scrolled right to column W and
# the active cell/selection sitting on Z14.
$win = $excel.ActiveWindow
$win.ScrollColumn = 23
$win.ScrollRow = 1
$ws.Range("Z14").Select()
